$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "La validación ... Hasta el momento se han validado ..." becomes
# "La validación ... Hasta el 16 de diciembre de 2021 se han validado ...".
# The target markup keeps the (unchanged) lead-in run, then places the new
# date and the trailing " se han validado " text into two brand-new runs,
# so we replace "momento se han validado " (everything from the word that
# must disappear through the end of the original run) with that exact
# two-run fragment via InsertXML, which lets us control run boundaries
# precisely instead of just swapping text inside the existing run.
$find = $d.Content.Find
$find.Execute("momento se han validado ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($find.Found) {
    $rng = $d.Range($find.Parent.Start, $find.Parent.End)

    $rPr = '<w:rPr><w:rFonts w:ascii="Source Sans Pro" w:hAnsi="Source Sans Pro" w:cs="Calibri"/><w:color w:val="6F7271"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-MX"/></w:rPr>'
    $runDate = '<w:r>' + $rPr + '<w:t>16 de diciembre de 2021</w:t></w:r>'
    $runTail = '<w:r>' + $rPr + '<w:t xml:space="preserve"> se han validado </w:t></w:r>'

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runDate + $runTail + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
}

# --- Change 2 -----------------------------------------------------------
# The run that hosts the inline picture gains <w:noProof/> in its run
# properties (keeps the spell/grammar checker from flagging the picture).
if ($d.InlineShapes.Count -ge 1) {
    $pic = $d.InlineShapes.Item(1)
    $pic.Range.NoProofing = 1
}
